$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1285.3077
$ws.Range("I121").Value = 1026.6666
$ws.Range("J121").Value = 1362.9
$ws.Range("K121").Value = 3079.9998
$ws.Range("L121").Value = 4088.7
$ws.Range("M121").Value = -1332.9998
$ws.Range("N121").Value = -7582.700000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3190.7908
$ws.Range("I132").Value = 2852.4666
$ws.Range("J132").Value = 3971.5386
$ws.Range("K132").Value = 8557.399800000001
$ws.Range("L132").Value = 11914.6158
$ws.Range("M132").Value = -6027.399800000001
$ws.Range("N132").Value = -16974.6158

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1435.766
$ws.Range("I137").Value = 1250.742
$ws.Range("J137").Value = 1794.25
$ws.Range("K137").Value = 3752.226
$ws.Range("L137").Value = 5382.75
$ws.Range("M137").Value = -1202.226
$ws.Range("N137").Value = -10482.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2648.8157
$ws.Range("I138").Value = 1719.8302
$ws.Range("J138").Value = 4789.522
$ws.Range("K138").Value = 5159.4906
$ws.Range("L138").Value = 14368.566
$ws.Range("M138").Value = -19.49060000000009
$ws.Range("N138").Value = -24648.566

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4822
$ws.Range("I141").Value = 2170.0857
$ws.Range("K141").Value = 6510.257100000001
$ws.Range("M141").Value = -1330.257100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11861.221
$ws.Range("I32").Value = 11817.057
$ws.Range("J32").Value = 12251.333
$ws.Range("K32").Value = 11817.057
$ws.Range("L32").Value = 12251.333
$ws.Range("M32").Value = -11530.057
$ws.Range("N32").Value = -12825.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 22225.4
$ws.Range("J112").Value = 22225.4
$ws.Range("L112").Value = 22225.4
$ws.Range("N112").Value = -25179.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H114").Value = 50000
$ws.Range("J114").Value = 50000
$ws.Range("L114").Value = 50000
$ws.Range("N114").Value = -58678

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 24915
$ws.Range("J124").Value = 24915
$ws.Range("L124").Value = 24915
$ws.Range("N124").Value = -34735

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2593.88
$ws.Range("I105").Value = 2739.7856
$ws.Range("J105").Value = 2408.182
$ws.Range("K105").Value = 2739.7856
$ws.Range("L105").Value = 2408.182
$ws.Range("M105").Value = -992.7856000000002
$ws.Range("N105").Value = -5902.182

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2971
$ws.Range("I134").Value = 3308.889
$ws.Range("J134").Value = 2780.9375
$ws.Range("K134").Value = 9926.667000000001
$ws.Range("L134").Value = 8342.8125
$ws.Range("M134").Value = -7391.667000000001
$ws.Range("N134").Value = -13412.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2359.9062
$ws.Range("I31").Value = 1565.762
$ws.Range("K31").Value = 1565.762
$ws.Range("M31").Value = -1270.762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2359.9062
$ws.Range("I34").Value = 1565.762
$ws.Range("K34").Value = 1565.762
$ws.Range("M34").Value = -1363.762

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 50000
$ws.Range("J111").Value = 50000
$ws.Range("L111").Value = 50000
$ws.Range("N111").Value = -58180

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1683.4706
$ws.Range("I134").Value = 1441.4584
$ws.Range("J134").Value = 2264.3
$ws.Range("K134").Value = 4324.3752
$ws.Range("L134").Value = 6792.900000000001
$ws.Range("M134").Value = -1789.3752
$ws.Range("N134").Value = -11862.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2312.7144
$ws.Range("I50").Value = 145.28572
$ws.Range("J50").Value = 4480.143
$ws.Range("K50").Value = 435.85716
$ws.Range("L50").Value = 13440.429
$ws.Range("M50").Value = 45.14283999999998
$ws.Range("N50").Value = -14402.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 2312.7144
$ws.Range("I53").Value = 145.28572
$ws.Range("J53").Value = 4480.143
$ws.Range("K53").Value = 435.85716
$ws.Range("L53").Value = 13440.429
$ws.Range("M53").Value = 45.14283999999998
$ws.Range("N53").Value = -14402.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1311.4445
$ws.Range("I121").Value = 1157.5
$ws.Range("J121").Value = 1434.6
$ws.Range("K121").Value = 3472.5
$ws.Range("L121").Value = 4303.799999999999
$ws.Range("M121").Value = -2162.5
$ws.Range("N121").Value = -6923.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7020.25
$ws.Range("I70").Value = 6537.8184
$ws.Range("J70").Value = 7609.8887
$ws.Range("K70").Value = 6537.8184
$ws.Range("L70").Value = 7609.8887
$ws.Range("M70").Value = -6267.8184
$ws.Range("N70").Value = -8149.8887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7020.25
$ws.Range("I73").Value = 6537.8184
$ws.Range("J73").Value = 7609.8887
$ws.Range("K73").Value = 6537.8184
$ws.Range("L73").Value = 7609.8887
$ws.Range("M73").Value = -5601.8184
$ws.Range("N73").Value = -9481.8887

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 18837.766
$ws.Range("J123").Value = 18837.766
$ws.Range("L123").Value = 18837.766
$ws.Range("N123").Value = -23737.766

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2493.037
$ws.Range("I132").Value = 1589.5714
$ws.Range("J132").Value = 3466
$ws.Range("K132").Value = 4768.7142
$ws.Range("L132").Value = 10398
$ws.Range("M132").Value = -2238.7142
$ws.Range("N132").Value = -15458

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 563.625
$ws.Range("I55").Value = 99.5
$ws.Range("J55").Value = 718.3333
$ws.Range("K55").Value = 99.5
$ws.Range("L55").Value = 718.3333
$ws.Range("M55").Value = 73.5
$ws.Range("N55").Value = -1064.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1436.9286
$ws.Range("I136").Value = 1265.2375
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 3795.7125
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -1245.7125
$ws.Range("N136").Value = -11700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 23113.25
$ws.Range("J64").Value = 23113.25
$ws.Range("L64").Value = 23113.25
$ws.Range("N64").Value = -23609.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 23113.25
$ws.Range("J67").Value = 23113.25
$ws.Range("L67").Value = 23113.25
$ws.Range("N67").Value = -24829.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22106.438
$ws.Range("J123").Value = 22106.438
$ws.Range("L123").Value = 22106.438
$ws.Range("N123").Value = -31906.438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 913.375
$ws.Range("I132").Value = 725.2162
$ws.Range("J132").Value = 1546.2727
$ws.Range("K132").Value = 2175.6486
$ws.Range("L132").Value = 4638.8181
$ws.Range("M132").Value = 354.3514
$ws.Range("N132").Value = -9698.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1234.6
$ws.Range("J136").Value = 800
$ws.Range("L136").Value = 2400
$ws.Range("N136").Value = -7500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 82150
$ws.Range("J138").Value = 82150
$ws.Range("L138").Value = 82150
$ws.Range("N138").Value = -92430
